$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Scattered F-column fixes (rows 19, 21, 23) ---
$ws.Range("F19").Value = 17.81
$ws.Range("F21").Value = ""
$ws.Range("F23").Value = 16.48

# --- Remove the "RM 232" row (row 26) entirely, shifting rows below up ---
$ws.Rows(26).Delete()

# --- Remove the "SC 92" row (now row 27 after the previous delete) ---
$ws.Rows(27).Delete()

# --- After the shifts, the row that now holds "SC 101" needs its F cell cleared ---
$ws.Range("F27").Value = ""

# --- And the last remaining row ("SC 232", now row 33) needs its F cell filled in ---
$ws.Range("F33").Value = 17.53
